$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# ALC row 33
$ws.Range("H33").Value = 33333596
$ws.Range("I33").Value = 33333596
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 33333596
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -33333367

# ALC row 46
$ws.Range("H46").Value = 1083.3334
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 1233.3334
$ws.Range("K46").Value = 2799.9999
$ws.Range("L46").Value = 3700.0002
$ws.Range("M46").Value = -2680.9999
$ws.Range("N46").Value = -3938.0002

# ALC row 60
$ws.Range("H60").Value = 1083.3334
$ws.Range("I60").Value = 933.3333
$ws.Range("J60").Value = 1233.3334
$ws.Range("K60").Value = 2799.9999
$ws.Range("L60").Value = 3700.0002
$ws.Range("M60").Value = -2315.9999
$ws.Range("N60").Value = -4668.0002

# ALC row 138
$ws.Range("H138").Value = 6100509
$ws.Range("I138").Value = 1482.5834
$ws.Range("J138").Value = 14710899
$ws.Range("K138").Value = 4447.7502
$ws.Range("L138").Value = 44132697
$ws.Range("M138").Value = 692.2497999999996
$ws.Range("N138").Value = -44142977

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 1915.8572
$ws.Range("I2").Value = 1551.8334
$ws.Range("J2").Value = 4100
$ws.Range("K2").Value = 1551.8334
$ws.Range("L2").Value = 4100
$ws.Range("M2").Value = -1438.8334
$ws.Range("N2").Value = -4326

# ARM row 16
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -213
$ws.Range("N16").Value = -1074

# ARM row 116
$ws.Range("H116").Value = 1915.8572
$ws.Range("I116").Value = 1551.8334
$ws.Range("J116").Value = 4100
$ws.Range("K116").Value = 1551.8334
$ws.Range("L116").Value = 4100
$ws.Range("M116").Value = 742.1666
$ws.Range("N116").Value = -8688

# ARM row 122
$ws.Range("H122").Value = 1889.5186
$ws.Range("I122").Value = 1888.5769
$ws.Range("J122").Value = 1914
$ws.Range("K122").Value = 5665.7307
$ws.Range("L122").Value = 5742
$ws.Range("M122").Value = -3215.7307
$ws.Range("N122").Value = -10642

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 1915.8572
$ws.Range("I3").Value = 1551.8334
$ws.Range("J3").Value = 4100
$ws.Range("K3").Value = 1551.8334
$ws.Range("L3").Value = 4100
$ws.Range("M3").Value = -1437.8334
$ws.Range("N3").Value = -4328

# BSM row 16
$ws.Range("H16").Value = 2169.6667
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 2203.6
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2203.6
$ws.Range("M16").Value = -1830
$ws.Range("N16").Value = -2543.6

# BSM row 20
$ws.Range("H20").Value = 992.125
$ws.Range("I20").Value = 1376.875
$ws.Range("J20").Value = 607.375
$ws.Range("K20").Value = 1376.875
$ws.Range("L20").Value = 607.375
$ws.Range("M20").Value = -1129.875
$ws.Range("N20").Value = -1101.375

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 707.375
$ws.Range("I22").Value = 431.8
$ws.Range("J22").Value = 1166.6666
$ws.Range("K22").Value = 431.8
$ws.Range("L22").Value = 1166.6666
$ws.Range("M22").Value = -81.80000000000001
$ws.Range("N22").Value = -1866.6666

# CRP row 132
$ws.Range("H132").Value = 3900.389
$ws.Range("I132").Value = 4039.25
$ws.Range("J132").Value = 3789.3
$ws.Range("K132").Value = 12117.75
$ws.Range("L132").Value = 11367.9
$ws.Range("M132").Value = -9587.75
$ws.Range("N132").Value = -16427.9

$ws = $wb.Worksheets.Item("CUL")
# CUL row 9
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -3448
$ws.Range("M9").ClearContents()

# CUL row 20
$ws.Range("H20").Value = 1080.5385
$ws.Range("I20").Value = 1019.8
$ws.Range("J20").Value = 1118.5
$ws.Range("K20").Value = 3059.4
$ws.Range("L20").Value = 3355.5
$ws.Range("M20").Value = -2832.4
$ws.Range("N20").Value = -3809.5

# CUL row 56
$ws.Range("H56").Value = 4050
$ws.Range("I56").Value = 4050
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4050
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -3520

# CUL row 70
$ws.Range("H70").Value = 3571.2354
$ws.Range("I70").Value = 2453
$ws.Range("J70").Value = 3915.3076
$ws.Range("K70").Value = 7359
$ws.Range("L70").Value = 11745.9228
$ws.Range("M70").Value = -7044
$ws.Range("N70").Value = -12375.9228

# CUL row 73
$ws.Range("H73").Value = 3571.2354
$ws.Range("I73").Value = 2453
$ws.Range("J73").Value = 3915.3076
$ws.Range("K73").Value = 7359
$ws.Range("L73").Value = 11745.9228
$ws.Range("M73").Value = -6267
$ws.Range("N73").Value = -13929.9228

$ws = $wb.Worksheets.Item("GSM")
# GSM row 14
$ws.Range("H14").Value = 4125501.2
$ws.Range("I14").Value = 7500002
$ws.Range("J14").Value = 751000.5
$ws.Range("K14").Value = 7500002
$ws.Range("L14").Value = 751000.5
$ws.Range("M14").Value = -7499834
$ws.Range("N14").Value = -751336.5

# GSM row 70
$ws.Range("H70").Value = 89485.14
$ws.Range("I70").Value = 122899.4
$ws.Range("J70").Value = 5949.5
$ws.Range("K70").Value = 122899.4
$ws.Range("L70").Value = 5949.5
$ws.Range("M70").Value = -122629.4
$ws.Range("N70").Value = -6489.5

# GSM row 73
$ws.Range("H73").Value = 89485.14
$ws.Range("I73").Value = 122899.4
$ws.Range("J73").Value = 5949.5
$ws.Range("K73").Value = 122899.4
$ws.Range("L73").Value = 5949.5
$ws.Range("M73").Value = -121963.4
$ws.Range("N73").Value = -7821.5

# GSM row 132
$ws.Range("H132").Value = 5589.864
$ws.Range("I132").Value = 5442.25
$ws.Range("J132").Value = 5674.2144
$ws.Range("K132").Value = 16326.75
$ws.Range("L132").Value = 17022.6432
$ws.Range("M132").Value = -13796.75
$ws.Range("N132").Value = -22082.6432

$ws = $wb.Worksheets.Item("LTW")
# LTW row 14
$ws.Range("H14").Value = 7500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 7500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 7500
$ws.Range("N14").Value = -7844

# LTW row 16
$ws.Range("H16").Value = 2161.04
$ws.Range("I16").Value = 1609.8422
$ws.Range("J16").Value = 3906.5
$ws.Range("K16").Value = 1609.8422
$ws.Range("L16").Value = 3906.5
$ws.Range("M16").Value = -1439.8422
$ws.Range("N16").Value = -4246.5

# LTW row 25
$ws.Range("H25").Value = 590.6667
$ws.Range("I25").Value = 386
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 386
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -156
$ws.Range("N25").Value = -1460

# LTW row 55
$ws.Range("H55").Value = 800
$ws.Range("I55").Value = 300
$ws.Range("J55").Value = 1050
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 1050
$ws.Range("M55").Value = -127
$ws.Range("N55").Value = -1396
